{"js": "// Apply the text replacements described by the diff.\n// Each pair is [old text, new text]; all old texts are unique within the\n// document, so a simple search+replace-text on each match is unambiguous.\nconst replacements = [\n  [\"2024-05-23 Thursday\", \"2024-05-24 Friday\"],\n  [\"25\u00f79=\", \"11\u00f79=\"],\n  [\"81\u00f78=\", \"40\u00f76=\"],\n  [\"89\u00f75=\", \"88\u00f77=\"],\n  [\"58\u00f75=\", \"90\u00f73=\"],\n  [\"69\u00f76=\", \"23\u00f78=\"],\n  [\"28\u00f77=\", \"12\u00f75=\"],\n  [\"25\u00f77=\", \"65\u00f75=\"],\n  [\"93\u00f72=\", \"77\u00f77=\"],\n  [\"73\u00f74=\", \"13\u00f77=\"],\n  [\"33\u00f78=\", \"27\u00f78=\"],\n  [\"91\u00f78=\", \"89\u00f73=\"],\n  [\"80\u00f76=\", \"60\u00f78=\"],\n  [\"11\u00f75=\", \"33\u00f74=\"],\n  [\"73\u00f76=\", \"20\u00f75=\"],\n  [\"43\u00f79=\", \"48\u00f73=\"],\n  [\"66\u00f75=\", \"64\u00f76=\"],\n  [\"37\u00f74=\", \"62\u00f77=\"],\n  [\"25\u00f74=\", \"78\u00f78=\"],\n  [\"63\u00f75=\", \"49\u00f77=\"],\n  [\"70\u00f72=\", \"55\u00f76=\"],\n  [\"66\u00f73=\", \"57\u00f72=\"],\n  [\"79\u00f77=\", \"13\u00f74=\"],\n  [\"38\u00f77=\", \"76\u00f74=\"],\n  [\"18\u00f74=\", \"60\u00f72=\"],\n  [\"21\u00f77=\", \"14\u00f77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff using Find/Replace.\n# Each pair is (old text, new text); all old texts are unique within the\n# document, so Find/Execute with ReplaceAll is unambiguous per pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-23 Thursday\", \"2024-05-24 Friday\"),\n    @(\"25\u00f79=\", \"11\u00f79=\"),\n    @(\"81\u00f78=\", \"40\u00f76=\"),\n    @(\"89\u00f75=\", \"88\u00f77=\"),\n    @(\"58\u00f75=\", \"90\u00f73=\"),\n    @(\"69\u00f76=\", \"23\u00f78=\"),\n    @(\"28\u00f77=\", \"12\u00f75=\"),\n    @(\"25\u00f77=\", \"65\u00f75=\"),\n    @(\"93\u00f72=\", \"77\u00f77=\"),\n    @(\"73\u00f74=\", \"13\u00f77=\"),\n    @(\"33\u00f78=\", \"27\u00f78=\"),\n    @(\"91\u00f78=\", \"89\u00f73=\"),\n    @(\"80\u00f76=\", \"60\u00f78=\"),\n    @(\"11\u00f75=\", \"33\u00f74=\"),\n    @(\"73\u00f76=\", \"20\u00f75=\"),\n    @(\"43\u00f79=\", \"48\u00f73=\"),\n    @(\"66\u00f75=\", \"64\u00f76=\"),\n    @(\"37\u00f74=\", \"62\u00f77=\"),\n    @(\"25\u00f74=\", \"78\u00f78=\"),\n    @(\"63\u00f75=\", \"49\u00f77=\"),\n    @(\"70\u00f72=\", \"55\u00f76=\"),\n    @(\"66\u00f73=\", \"57\u00f72=\"),\n    @(\"79\u00f77=\", \"13\u00f74=\"),\n    @(\"38\u00f77=\", \"76\u00f74=\"),\n    @(\"18\u00f74=\", \"60\u00f72=\"),\n    @(\"21\u00f77=\", \"14\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
